# Guide.docx update — add missing end-of-sentence periods and two new
# explanatory paragraphs (Vietnamese), per commit "Update Guide version 2".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Heading "Cài đặt và tạo cơ sở dữ liệu" -> add trailing period, then
#    insert a new (non-bold) explanatory paragraph right after it.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(2)
$r1 = $p1.Range
$r1.MoveEnd(1, -1)            # exclude the paragraph mark
$r1.Collapse(0)
$r1.InsertAfter(".")
$r1.Font.Bold = 1
$r1.Font.BoldBi = 1

# Insert the new paragraph after paragraph 2.
$r1.Collapse(0)
$r1.InsertParagraphAfter()
$newPara1 = $p1.Next()
$newRange1 = $newPara1.Range
$newRange1.ParagraphFormat.Style = "ListParagraph"
$newRange1.ParagraphFormat.LeftIndent = 36
$newRange1.ParagraphFormat.FirstLineIndent = 0
$newRange1.ParagraphFormat.Alignment = 0
$newRange1.MoveEnd(1, -1)
$newRange1.InsertAfter("Để để bắt đầu làm việc bạn cần cài đặt chương trình tạo server và quản lý cơ sở dữ liệu.")
$newRange1.Font.Bold = 0
$newRange1.Font.BoldBi = 0

# ---------------------------------------------------------------------
# 2) "Tải composer" -> add trailing period (bold run, matches heading).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Tải composer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(".")
$rng.Font.Bold = 1
$rng.Font.BoldBi = 1

# ---------------------------------------------------------------------
# 3) "Sau đó cài đặt" -> add trailing period (regular run).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Sau đó cài đặt", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(".")

# ---------------------------------------------------------------------
# 4) "composer –v" -> add trailing period (regular run).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("composer –v", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(".")

# ---------------------------------------------------------------------
# 5) "... artisan migrate" -> add trailing period (regular run).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("artisan migrate", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(".")

# ---------------------------------------------------------------------
# 6) After heading "4.  Tạo tài khoản quản lý", insert a new paragraph
#    (Normal style, first-line indent, justified) before "Bước 1: ...".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Tạo tài khoản quản lý", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingPara = $rng.Paragraphs(1)
$hr = $headingPara.Range
$hr.Collapse(0)
$hr.InsertParagraphAfter()
$newPara2 = $headingPara.Next()
$newRange2 = $newPara2.Range
$newRange2.ParagraphFormat.Style = "Normal"
$newRange2.ParagraphFormat.RightIndent = 0
$newRange2.ParagraphFormat.FirstLineIndent = 36
$newRange2.ParagraphFormat.Alignment = 3
$newRange2.MoveEnd(1, -1)
$newRange2.InsertAfter("Sau khi thiết lập hệ thống, chúng ta có thể bắt đầu làm việc với web sau trong phần mềm trình duyệt.")
$newRange2.Font.Bold = 0
$newRange2.Font.BoldBi = 0
$newRange2.Collapse(0)
$newRange2.InsertAfter("`t")
